# CFWH and OFWH Efficiencies
# Fill in the 2nd-law component efficiency values (and the still-missing
# "N/A" / unit placeholders) on the "Grade C cover page" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Grade C cover page")
$ws.Activate()

# --- Closed FWH (CFWH) component efficiencies: rows 25-30 ---
# condenser (row 28) and superheater (row 30) are "N/A"
$ws.Range("C28").Value = "N/A"
$ws.Range("C30").Value = "N/A"

# --- Open FWH (OFWH) section: rows 33-44 ---
$ws.Range("C33").Value = 7.6      # recommended OFWH pressure
$ws.Range("C34").Value = 22.6     # open FWH bleedoff ratio, y
$ws.Range("C35").Value = 37.3     # thermal cycle efficiency
$ws.Range("C36").Value = 69.3     # exergetic cycle efficiency

# 2nd law component efficiencies (OFWH): rows 38-44
$ws.Range("C38").Value = 3.8      # turbine 1
$ws.Range("C39").Value = 5.5      # feedwater heater
$ws.Range("C40").Value = 90.5     # turbine 2
$ws.Range("C41").Value = "N/A"    # condenser
$ws.Range("C42").Value = 95.4     # pump 1
$ws.Range("C43").Value = 96.6     # pump 2
$ws.Range("D43").Value = "%"      # previously-missing unit label
$ws.Range("C44").Value = "N/A"    # superheater

# Restore the selection/active cell to match the finished state of the edit
$ws.Range("C44").Select()
